# Updated cryptos list on Sun Jun 11 11:53:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.825.32"

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.754.15"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'237.02"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.05%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.5097"
$ws.Range("E7").Value = "  +3.16%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2700"
$ws.Range("E8").Value = "  +8.59%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06204"
$ws.Range("E9").Value = "  +3.75%  "

# Row 10 - WrappedEther
$ws.Range("D10").Value = "1.745.88"
$ws.Range("E10").Value = "  -0.08%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.06923"
$ws.Range("E11").Value = "  +2.12%  "

# Row 12 - Solana
$ws.Range("E12").Value = "  +4.91%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "'0.6295"
$ws.Range("E13").Value = "  +8.05%  "

# Row 14 & 15 - Litecoin / Polkadot swapped places
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.499"
$ws.Range("E14").Value = "  +0.51%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'78.43"
$ws.Range("E15").Value = "  +1.55%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  +0.02%  "

# Row 17 - Dai
$ws.Range("E17").Value = "  +0.07%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.845.79"

# Row 19 - Avalanche
$ws.Range("E19").Value = "  +1.03%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.000006725"
$ws.Range("E20").Value = "  +2.40%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "1.969.25"
$ws.Range("E21").Value = "  +0.18%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.079"
$ws.Range("E22").Value = "  +2.26%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "'8.271"
$ws.Range("E23").Value = "  +4.36%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'5.185"
$ws.Range("E24").Value = "  +2.81%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'136.65"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - EthereumClassic
$ws.Range("E26").Value = "  +5.37%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -1.84%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "'1.791"
$ws.Range("E28").Value = "  -2.49%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "'102.79"
$ws.Range("E29").Value = "  +1.75%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "'0.08273"
$ws.Range("E30").Value = "  +1.87%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'3.732"
$ws.Range("E31").Value = "  -1.81%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.431"
$ws.Range("E32").Value = "  +2.32%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.04415"
$ws.Range("E33").Value = "  -0.12%  "

# Row 34 - Frax
$ws.Range("D34").Value = "'0.9991"
$ws.Range("E34").Value = "  +0.03%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'2.649"
$ws.Range("E35").Value = "  -0.57%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "'1.006"
$ws.Range("E36").Value = "  -1.13%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'0.6057"
$ws.Range("E37").Value = "  -0.30%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "'2.695"
$ws.Range("E38").Value = "  -0.13%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'1.971"
$ws.Range("E39").Value = "  -4.34%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "'0.01560"
$ws.Range("E40").Value = "  +3.90%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.13%  "

# Row 42 - Quant
$ws.Range("D42").Value = "'102.26"
$ws.Range("E42").Value = "  -1.50%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "'0.3867"
$ws.Range("E43").Value = "  +2.52%  "

# Row 44 - TrustWalletToken
$ws.Range("D44").Value = "'0.7501"
$ws.Range("E44").Value = "  -3.35%  "

# Row 45 - FraxShare
$ws.Range("D45").Value = "'4.904"
$ws.Range("E45").Value = "  -5.46%  "

# Row 46 - Cronos
$ws.Range("D46").Value = "'0.05512"
$ws.Range("E46").Value = "  +7.52%  "

# Row 47 - Algorand
$ws.Range("D47").Value = "'0.1097"
$ws.Range("E47").Value = "  +1.41%  "

# Row 48 - Aptos
$ws.Range("D48").Value = "'5.979"
$ws.Range("E48").Value = "  -0.07%  "

# Row 49 - Elrond
$ws.Range("D49").Value = "'30.26"
$ws.Range("E49").Value = "  -0.54%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'52.90"
$ws.Range("E50").Value = "  +0.44%  "
